$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest crypto price/volume snapshot: cell -> new value.
# Cells are forced to Text number format before the assignment so
# numeric-looking strings (e.g. "574.94") are stored as literal text
# instead of being reinterpreted as numbers, matching the source data.
# The style is then reset to Normal so no stray formatting is left on the cell.
$updates = @{
    'D2' = '61.466.49'
    'E2' = '  +0.74%  '
    'D3' = '3.384.17'
    'E3' = '  -0.21%  '
    'E4' = '  -0.01%  '
    'D5' = '574.94'
    'E5' = '  +0.46%  '
    'D6' = '140.90'
    'E6' = '  -0.94%  '
    'E7' = '  +0.03%  '
    'D8' = '0.473'
    'E8' = '  -0.52%  '
    'E9' = '  +0.48%  '
    'E10' = '  -1.18%  '
    'D11' = '0.385'
    'E11' = '  -2.71%  '
    'D12' = '3.962.45'
    'E12' = '  -0.23%  '
    'E13' = '  +0.08%  '
    'D14' = '28.15'
    'E14' = '  +1.17%  '
    'D15' = '3.365.04'
    'E15' = '  -0.94%  '
    'E16' = '  -0.98%  '
    'D17' = '61.528.47'
    'E17' = '  +0.83%  '
    'D18' = '6.11'
    'E18' = '  -0.38%  '
    'E19' = '  -1.45%  '
    'E20' = '  +0.50%  '
    'D21' = '389.27'
    'E21' = '  +1.78%  '
    'D22' = '74.91'
    'E22' = '  +0.63%  '
    'E23' = '  -1.49%  '
    'E24' = '  +0.25%  '
    'E25' = '  +9.02%  '
    'E26' = '  -4.17%  '
    'E27' = '  +0.30%  '
    'E28' = '  -0.81%  '
    'E29' = '  -0.16%  '
    'E30' = '  -0.70%  '
    'E31' = '  -0.19%  '
    'D32' = '0.999'
    'E32' = '  +0.01%  '
    'E33' = '  -1.09%  '
    'E34' = '  -1.52%  '
    'D35' = '168.59'
    'E35' = '  +0.91%  '
    'E36' = '  +0.11%  '
    'D37' = '3.416.92'
    'E37' = '  -0.19%  '
    'D38' = '1.46'
    'E38' = '  -1.30%  '
    'E39' = '  -0.89%  '
    'D40' = '26.03'
    'E40' = '  -4.60%  '
    'D41' = '0.777'
    'E41' = '  -0.47%  '
    'D42' = '4.43'
    'E42' = '  +0.21%  '
    'E43' = '  -1.73%  '
    'E44' = '  +1.22%  '
    'D45' = '2.463.41'
    'E45' = '  -0.89%  '
    'D46' = '22.59'
    'E46' = '  -1.26%  '
    'E47' = '  -2.37%  '
    'E48' = '  +0.03%  '
    'E49' = '  -1.04%  '
    'D50' = '2.01'
    'E50' = '  -5.27%  '
    'E51' = '  -1.79%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
